# "Seventeenth Commit: Completed my rectangle test plan."
# Fill in the Developer field and the Method Inputs / Actual Result columns
# of the Rectangle unit-test plan, then leave the selection where the
# author finished working (G13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Developer field -------------------------------------------------
$ws.Range("C3").Value = "Michael Obikwere"

# --- Test rows 7-13 (Method Inputs = E, Actual Result = F, Expected = G) --
# Values are entered in the same order the shared-string table records
# them so duplicate strings resolve to the same entry.

# Row 7 - __init__ / attribute set test
$ws.Range("E7").Value = "None"
$ws.Range("E11").Value = "Rectangle(""Red"", 8, 10)"
$ws.Range("G8").Value = "ValueError"
$ws.Range("G7").Value = "Attributes Set"
$ws.Range("F7").Value = "color : ""Red""`nlength : 8`nwidth : 10"
$ws.Range("F10").Value = "Color : ""Red""`nlength : 8`nwidth : ""10"""
$ws.Range("F9").Value = "color : ""Red""`nlength : ""8""`nwidth : 10"
$ws.Range("F8").Value = "color : """"`nlength : 8`nwidth : 10"
$ws.Range("G11").Value = "The shape color is red.`nThis rectangle has four sides with the lengths of 8, 10, 8 and 10 centimeters."
$ws.Range("G12").Value = "area = 80"
$ws.Range("G13").Value = "perimeter = 36"

# Remaining duplicate cells (same shared strings as above).
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("F11").Value = "None"
$ws.Range("F12").Value = "None"
$ws.Range("F13").Value = "None"

$ws.Range("E12").Value = "Rectangle(""Red"", 8, 10)"
$ws.Range("E13").Value = "Rectangle(""Red"", 8, 10)"

$ws.Range("G9").Value = "ValueError"
$ws.Range("G10").Value = "ValueError"

# Row 13 picks up the bold "filled-in" look the other test rows already
# had (font only differs by Bold between style 3 and style 9).
$ws.Range("E13:G13").Font.Bold = $true

# --- Row heights (autofit-style tweaks from the re-save) -------------
# Rows 1/3/4 only carry a thick-bottom border (no real customHeight in
# either revision) so they are left alone; row 2 and 13-23 already had
# an explicit customHeight and just shrink slightly.
$ws.Rows.Item(2).RowHeight = 73.2
$ws.Rows.Item(13).RowHeight = 31.2
$ws.Rows.Item(14).RowHeight = 31.2
$ws.Rows.Item(15).RowHeight = 31.2
$ws.Rows.Item(16).RowHeight = 31.2
$ws.Rows.Item(17).RowHeight = 31.2
$ws.Rows.Item(18).RowHeight = 31.2
$ws.Rows.Item(19).RowHeight = 31.2
$ws.Rows.Item(20).RowHeight = 31.2
$ws.Rows.Item(21).RowHeight = 31.2
$ws.Rows.Item(22).RowHeight = 31.2
$ws.Rows.Item(23).RowHeight = 31.2

# --- Leave the selection where the author finished --------------------
$ws.Range("G13").Select()
